$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = 0.03858115821364612
$ws.Range("J18").Value = 0.23502451085558
$ws.Range("K18").Value = -0.2184581852699795
$ws.Range("L18").Value = 2.731007526202857

$ws.Range("I19").Value = 0.3569533601582136
$ws.Range("J19").Value = 0.503493624976924
$ws.Range("K19").Value = 0.1370214183447129
$ws.Range("L19").Value = 2.132386593961161
